# Replace the benchmark results with a new calculation run and a new
# method ("kmeans SOM" is recomputed, "random kmeans", "kde kmeans" and
# "kde kmeans SOM" are dropped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (B1:G1) ------------------------------------------------
$headers = @("kmeans++", "random SOM", "kmeans SOM", "kmeans++ SOM", "SOM++", "kde SOM")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, 2 + $c).Value = $headers[$c]
}

# --- New data rows (B2:G9) --------------------------------------------------
$rows = @{
    2 = @(50, 50, 50, 50, 50, 50)
    3 = @(0.2934, 0.2847, -1, 0.4207, 0.3586, 0.4707)
    4 = @(0.0598, 0.2466, 0, 0.3431, 0.143, 0.0985)
    5 = @(0.2304, -0.3911, -1, -1, 0.1212, 0.2442)
    6 = @(0.2639, 0.1191, -1, 0.2919, 0.2621, 0.397)
    7 = @(0.2772, 0.3326, -1, 0.4218, 0.3472, 0.4317)
    8 = @(0.2907, 0.42, -1, 0.6521, 0.4283, 0.5809)
    9 = @(0.471, 0.7929, -1, 0.8189, 0.8189, 0.611)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, 2 + $c).Value = $vals[$c]
    }
}

# --- Drop the old H:J columns (random kmeans / kde kmeans / kde kmeans SOM
#     no longer exist as separate methods) so the used range shrinks back
#     down to A1:G9. -------------------------------------------------------
$ws.Range("H1:J9").Clear()
